# highPearson.xlsx — refresh the correlation sample with a re-run of the
# underlying analysis: updated Pearson correlation coefficients for the
# existing keyword rows, one keyword swapped ("estou com covid" ->
# "peguei covid"), a couple of low-signal keyword rows trimmed off the
# bottom, and the now-empty D12 score cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: teste covid ---
$ws.Range("C2").Value = 0.7875854822277522
$ws.Range("D2").Value = 0.7975381766671763

# --- Row 3: exame covid ---
$ws.Range("C3").Value = 0.7473921235971395
$ws.Range("D3").Value = 0.876625308755157
$ws.Range("E3").Value = 0.7684654722085865

# --- Row 4: teste rápido covid ---
$ws.Range("C4").Value = 0.8163223847304697

# --- Row 5: igg ---
$ws.Range("C5").Value = 0.9019911760830701
$ws.Range("D5").Value = 0.765061779523126

# --- Row 6: igm ---
$ws.Range("C6").Value = 0.9038455621815769
$ws.Range("D6").Value = 0.7735404525439052

# --- Row 7: igg igm ---
$ws.Range("C7").Value = 0.9066945790484944
$ws.Range("D7").Value = 0.7593155893104981

# --- Row 8: covid igg ---
$ws.Range("C8").Value = 0.8750976242612639
$ws.Range("D8").Value = 0.7369821730538669

# --- Row 9: covid igm igg ---
$ws.Range("C9").Value = 0.8826151731133121
$ws.Range("D9").Value = 0.7542000109207299

# --- Row 10: reagente igg ---
$ws.Range("C10").Value = 0.8394037122836262
$ws.Range("D10").Value = 0.8945000158630423
$ws.Range("E10").Value = 0.7429684860210746

# --- Row 11: teste igg ---
$ws.Range("C11").Value = 0.8310669428757275

# --- Row 12: exame igg (D12 score no longer present) ---
$ws.Range("C12").Value = 0.8286654517918608
$ws.Range("D12").Value = ""

# --- Row 13: exame cotonete ---
$ws.Range("D13").Value = 0.8119023202767103
$ws.Range("E13").Value = 0.8514641675341623

# --- Row 14: covid pcr ---
$ws.Range("D14").Value = 0.8560927593465805
$ws.Range("E14").Value = 0.8460074840257149

# --- Row 15: pcr exame covid ---
$ws.Range("D15").Value = 0.8964147451952494
$ws.Range("E15").Value = 0.8450725592291811

# --- Row 16: pcr ---
$ws.Range("D16").Value = 0.8343129232104819
$ws.Range("E16").Value = 0.8440863383087603

# --- Row 17: keyword swapped, "estou com covid" -> "peguei covid"; gains E17 ---
$ws.Range("B17").Value = "peguei covid"
$ws.Range("D17").Value = 0.7523444142781768
$ws.Range("E17").Value = 0.7133486891517318

# --- Row 18: "pico brasil" -> "brasil coronavírus" ---
$ws.Range("B18").Value = "brasil coronavírus"
$ws.Range("E18").Value = -0.7101778375410028

# --- Row 19: "pico coronavírus" -> "brasil coronavírus mortes" ---
$ws.Range("B19").Value = "brasil coronavírus mortes"
$ws.Range("E19").Value = -0.7005033557810264

# --- Rows 20-25 (the remaining, lower-correlation keyword rows) are dropped ---
$ws.Range("A20:F25").Delete()
